$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.188.89"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.793.76"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'327.34"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "'0.4516"
$ws.Range("E7").Value = "  +16.55%  "
$ws.Range("D8").Value = "'0.3739"
$ws.Range("E8").Value = "  +10.11%  "
$ws.Range("D9").Value = "'44.72"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'1.143"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'0.07548"
$ws.Range("E11").Value = "  +4.69%  "
$ws.Range("D12").Value = "'22.48"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "'6.286"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "'7.545"
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("D16").Value = "1.792.87"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "'0.00001089"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "'0.06745"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "'80.86"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'17.49"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Value = "'6.345"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "28.208.27"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'11.78"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'2.425"
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Value = "'20.52"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("D27").Value = "'151.68"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "'2.350"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").Value = "1.996.67"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("D30").Value = "'132.66"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("D31").Value = "'1.233"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "'4.021"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'0.09412"
$ws.Range("E33").Value = "  +8.07%  "
$ws.Range("D34").Value = "'5.788"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").Value = "'0.2353"
$ws.Range("E35").Value = "  +11.64%  "
$ws.Range("D36").Value = "'12.11"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'0.06319"
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("D38").Value = "'0.02326"
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").Value = "'5.173"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'0.6560"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "'8.335"
$ws.Range("E41").Value = "  +6.17%  "
$ws.Range("D42").Value = "'1.477"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'1.201"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'14.14"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "'0.6092"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").Value = "'3.785"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'129.70"
$ws.Range("E48").Value = "  +2.42%  "
$ws.Range("D49").Value = "'2.023"
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("D50").Value = "'0.07127"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "'1.159"
$ws.Range("E51").Value = "  +0.36%  "
